# Intensiv.xlsx — add a new "10_11_2020" column of admission figures
# right after the "Aldersgruppe" label column (new column B), shifting
# all the existing date columns one place to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column B ("17_11_2020").
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 2).Value = "10_11_2020"

# New column's data values (admissions on 10_11_2020), row by row.
$newColumnValues = @(0, 1, 9, 15, 33, 92, 132, 191, 67, 2)

for ($i = 0; $i -lt $newColumnValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $newColumnValues[$i]
}

# Excel extends the "I alt" total-row formula pattern into the freshly
# inserted column, same as it does for the other date columns.
$ws.Cells.Item(12, 2).Formula = "=SUM(B2:B11)"

# Mirror Excel's default column-insert behaviour, which copies
# formatting from the column to the left (column A) into the freshly
# inserted column — this is what gives the new data cells the same
# quote-prefix style as the "Aldersgruppe" labels. Done last so the
# value-assignments above don't clobber the copied format.
$ws.Range("A2:A11").Copy()
$ws.Range("B2:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The sheet view no longer pins a fixed top-left scroll cell, and the
# selection moves to C18.
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C18").Select()
